# =====================================================================
# Apply the NCPI CodeSystem spreadsheet update:
#   * Metadata sheet: refresh URL + Date values
#   * Concepts sheet: reword 3 existing definitions and add 5 new
#     concept rows (SubStudyCount, VariableCount, AnalysesCount,
#     MolecularDatasetCount, PhenotypeDatasetCount)
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata": update the URL and Date property values
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Cells.Item(2, 2).Value = "https://nih-ncpi.github.io/ncpi-fhir-ig/CodeSystem/ncpi"
$ws1.Cells.Item(8, 2).Value = "2022-10-31T22:46:37+00:00"

# ---------------------------------------------------------------------
# Sheet "Concepts": reword definitions + insert new concept rows
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Concepts")

function Insert-ConceptRows {
    param($AfterRow, $HowMany, $LastDataRow)
    # Push existing rows (AfterRow+1 .. LastDataRow) down by $HowMany,
    # working bottom-to-top so source rows are never clobbered before
    # they're read. Uses Copy + PasteSpecial(Values) + PasteSpecial(Formats)
    # instead of Rows.Insert() so no new cell style gets fabricated
    # (matches the fact that styles.xml is untouched by this edit).
    for ($r = $LastDataRow; $r -ge ($AfterRow + 1); $r--) {
        $srcAddr = "A" + $r + ":D" + $r
        $dstAddr = "A" + ($r + $HowMany) + ":D" + ($r + $HowMany)
        $ws2.Range($srcAddr).Copy()
        $ws2.Range($dstAddr).PasteSpecial(-4163)
        $ws2.Range($srcAddr).Copy()
        $ws2.Range($dstAddr).PasteSpecial(-4122)
    }
}

# Reword the existing CohortCount / SampleCount / Participant definitions
$ws2.Cells.Item(6, 3).Value = "Number of cohorts that are part of this study"
$ws2.Cells.Item(7, 3).Value = "Number of samples that are part of this study"
$ws2.Cells.Item(8, 3).Value = "Number of participants that are part of this study"

# Insert a new row for SubStudyCount right after CohortCount (row 6)
Insert-ConceptRows 6 1 14
$ws2.Cells.Item(7, 1).Value = "1"
$ws2.Cells.Item(7, 2).Value = "SubStudyCount"
$ws2.Cells.Item(7, 3).Value = "Number of sub-studies that are part of this study"
$ws2.Cells.Item(7, 4).Value = ""

# Insert four new rows after Participant (now row 9, since the table grew by
# one row above): VariableCount, AnalysesCount, MolecularDatasetCount,
# PhenotypeDatasetCount
Insert-ConceptRows 9 4 15

$ws2.Cells.Item(10, 1).Value = "1"
$ws2.Cells.Item(10, 2).Value = "VariableCount"
$ws2.Cells.Item(10, 3).Value = "Number of variables that are part of this study"
$ws2.Cells.Item(10, 4).Value = ""

$ws2.Cells.Item(11, 1).Value = "1"
$ws2.Cells.Item(11, 2).Value = "AnalysesCount"
$ws2.Cells.Item(11, 3).Value = "Number of analyses that are part of this study"
$ws2.Cells.Item(11, 4).Value = ""

$ws2.Cells.Item(12, 1).Value = "1"
$ws2.Cells.Item(12, 2).Value = "MolecularDatasetCount"
$ws2.Cells.Item(12, 3).Value = "Number of molecular datasets that are part of this study"
$ws2.Cells.Item(12, 4).Value = ""

$ws2.Cells.Item(13, 1).Value = "1"
$ws2.Cells.Item(13, 2).Value = "PhenotypeDatasetCount"
$ws2.Cells.Item(13, 3).Value = "Number of phenotype datasets that are part of this study"
$ws2.Cells.Item(13, 4).Value = ""
